$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataCombined")
$ws.Activate()
$ws.Range("E2").Value = "Organism|PeripheralVenousBlood|Aciclovir|Plasma (Peripheral Venous Blood)"
$ws.Range("E3").Select()
